# Applies updated simulation results (Nausicaa v3.0 glider design) to column B values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = [double]"4.586677829652073"
$ws.Range("B3").Value = [double]"8.000000079969384"
$ws.Range("B4").Value = [double]"65.00000064970294"
$ws.Range("B6").Value = [double]"2.366201640692703"
$ws.Range("B7").Value = [double]"0.1302879996247762"
$ws.Range("B8").Value = [double]"12.01608239753199"
$ws.Range("B9").Value = [double]"0.9032065732239588"
$ws.Range("B10").Value = [double]"0.03999999003510714"
$ws.Range("B11").Value = [double]"0.3999999900893299"
$ws.Range("B12").Value = [double]"0.019999990001299"
$ws.Range("B13").Value = [double]"47617.94750078177"
$ws.Range("B14").Value = [double]"-0.1351272611301569"
$ws.Range("B15").Value = [double]"0.7187750035881874"
$ws.Range("B16").Value = [double]"1.572085127265426"
$ws.Range("B17").Value = [double]"1.244446304013794"
$ws.Range("B18").Value = [double]"0.1516489098467315"
$ws.Range("B20").Value = [double]"0.2384049957360621"
$ws.Range("B21").Value = [double]"0.7187750035881874"
$ws.Range("B22").Value = [double]"0.283687908460051"
$ws.Range("B23").Value = [double]"0.07092197711501275"
$ws.Range("B25").Value = [double]"0.02011970735160956"
$ws.Range("B26").Value = [double]"0.7187750035881874"
$ws.Range("B27").Value = [double]"0.144420588295566"
$ws.Range("B28").Value = [double]"0.07221029414778299"
$ws.Range("B30").Value = [double]"0.01042865316181868"
$ws.Range("B31").Value = [double]"-0.6788869189165445"
$ws.Range("B32").Value = [double]"5.901049369916801e-07"
$ws.Range("B33").Value = [double]"1.379937337878626e-06"
$ws.Range("B34").Value = [double]"3.164437083943233"
$ws.Range("B39").Value = [double]"1.589021735948907e-09"
$ws.Range("B40").Value = [double]"0.03683269640263574"
$ws.Range("B41").Value = [double]"-1.795779200810796e-17"
$ws.Range("B42").Value = [double]"0.006065873232640401"
$ws.Range("B43").Value = [double]"0.01648895843197087"
$ws.Range("B44").Value = [double]"0.003226500725865624"
$ws.Range("B45").Value = [double]"0.01968929380950151"
$ws.Range("B46").Value = [double]"7.54530294616237e-22"
$ws.Range("B47").Value = [double]"-6.933513753893814e-05"
$ws.Range("B48").Value = [double]"2.410209944985826e-20"
$ws.Range("B49").Value = [double]"0.07420932212624523"
$ws.Range("B50").Value = [double]"0.002580667540697819"
$ws.Range("B51").Value = [double]"0.001361936639877316"
$ws.Range("B56").Value = [double]"0.007685120382465098"
$ws.Range("B58").Value = [double]"-9.999677911335569e-09"
$ws.Range("B59").Value = [double]"0.009650962935168605"
$ws.Range("B60").Value = [double]"-0.03791222746168288"
$ws.Range("B63").Value = [double]"-0.03791222746168288"
$ws.Range("B64").Value = [double]"0.7858571652018056"
$ws.Range("B65").Value = [double]"0.01707125488268932"
$ws.Range("B66").Value = [double]"0.6478530264731746"
$ws.Range("B69").Value = [double]"0.6478530264731746"
$ws.Range("B70").Value = [double]"0.1418439542300255"
$ws.Range("B72").Value = [double]"0.6465647094404043"
$ws.Range("B75").Value = [double]"0.6465647094404043"
$ws.Range("B77").Value = [double]"0.144420588295566"
$ws.Range("B78").Value = [double]"-0.1716630531957717"
$ws.Range("B79").Value = [double]"-6.526460340647266e-09"
$ws.Range("B80").Value = [double]"3.029898005867883"
$ws.Range("B81").Value = [double]"3.026881886503915e-08"
$ws.Range("B82").Value = [double]"1.500590016828562e-15"
$ws.Range("B83").Value = [double]"7.967749093537576e-10"
$ws.Range("B84").Value = [double]"0.1716630531957717"
$ws.Range("B85").Value = [double]"-0.251687865128"
$ws.Range("B86").Value = [double]"-3.026881886503915e-08"
$ws.Range("B87").Value = [double]"-3.008513443069427e-08"
$ws.Range("B88").Value = [double]"3.024302125836964"
$ws.Range("B89").Value = [double]"-6.526460340647266e-09"
$ws.Range("B90").Value = [double]"0.251687865128"
$ws.Range("B91").Value = [double]"-3.026881886503915e-08"
$ws.Range("B92").Value = [double]"1.500590016828562e-15"
$ws.Range("B93").Value = [double]"-7.967749093537576e-10"
$ws.Range("B94").Value = [double]"0.9844814942268494"
$ws.Range("B95").Value = [double]"-2.12451638785743e-09"
$ws.Range("B96").Value = [double]"0.081930321518855"
$ws.Range("B97").Value = [double]"-6.267607586438256e-09"
$ws.Range("B98").Value = [double]"3.22110693120042e-15"
$ws.Range("B99").Value = [double]"-1.649840546740741e-10"
$ws.Range("B100").Value = [double]"2.893252895166951"
$ws.Range("B101").Value = [double]"0.1194706972139097"
$ws.Range("B102").Value = [double]"-1.053337804132748e-09"
$ws.Range("B103").Value = [double]"-2.991480078928799e-08"
$ws.Range("B104").Value = [double]"0.08204659517439658"
$ws.Range("B105").Value = [double]"-4.150885159898143e-09"
$ws.Range("B106").Value = [double]"1.571714374132259"
$ws.Range("B107").Value = [double]"0.7730833692395932"
$ws.Range("B108").Value = [double]"0.1300667733447216"
$ws.Range("B109").Value = [double]"0.009964522977254761"
$ws.Range("B111").Value = [double]"-1.951563910473908e-18"
$ws.Range("B112").Value = [double]"-0.08171383071985455"
$ws.Range("B113").Value = [double]"-1.592421940838085e-19"
$ws.Range("B114").Value = [double]"0.283687908460051"
$ws.Range("B117").Value = [double]"0.002912045769335567"
$ws.Range("B118").Value = [double]"-5.473122536514518e-09"
$ws.Range("B119").Value = [double]"-3.54018073779551e-10"
$ws.Range("B120").Value = [double]"-6.184121745671928e-05"
$ws.Range("B121").Value = [double]"3.354110250687772e-09"
$ws.Range("B122").Value = [double]"0.144420588295566"
$ws.Range("B124").Value = [double]"0.0009824573252918201"
$ws.Range("B125").Value = [double]"0.001029734088083553"
$ws.Range("B127").Value = [double]"-2.004462042719001e-20"
$ws.Range("B128").Value = [double]"-0.0002709232370838084"
$ws.Range("B129").Value = [double]"1.585639188626297e-20"
$ws.Range("B132").Value = [double]"0.1333770000485835"
$ws.Range("B133").Value = [double]"0.1183108650794165"
$ws.Range("B134").Value = [double]"5.668814753527495"
$ws.Range("B135").Value = [double]"0.4869366765324509"
$ws.Range("B136").Value = [double]"1.72826863030896e-09"
$ws.Range("B137").Value = [double]"3.306101893368849e-08"
$ws.Range("B138").Value = [double]"-0.2267525336519681"
$ws.Range("B139").Value = [double]"-2.499074371224133e-09"
$ws.Range("B140").Value = [double]"0.04289865128533987"
$ws.Range("B141").Value = [double]"-6.613047460355122e-06"
$ws.Range("B142").Value = [double]"-1.47418709485834e-07"
$ws.Range("B143").Value = [double]"0.07498834044235578"
$ws.Range("B144").Value = [double]"-0.02499999012464939"
$ws.Range("B145").Value = [double]"7.935290679272146e-07"
$ws.Range("B146").Value = [double]"0.0155638692857712"
$ws.Range("B147").Value = [double]"-0.2894543941307604"
$ws.Range("B148").Value = [double]"-0.003999685263766927"
$ws.Range("B149").Value = [double]"-0.00216724750756736"
$ws.Range("B150").Value = [double]"-0.06588786792411054"
$ws.Range("B151").Value = [double]"-0.8845379351899274"
$ws.Range("B152").Value = [double]"4.34759541575086e-05"
$ws.Range("B153").Value = [double]"-0.2759381839222994"
$ws.Range("B154").Value = [double]"0.546945228151907"
$ws.Range("B155").Value = [double]"0.3665173469717813"
$ws.Range("B156").Value = [double]"-1.823992013740275e-09"
$ws.Range("B157").Value = [double]"-1.345152199905187e-08"
$ws.Range("B158").Value = [double]"-14.4710600542043"
$ws.Range("B159").Value = [double]"2.037423588680084e-09"
$ws.Range("B160").Value = [double]"-1.463320475902918e-05"
$ws.Range("B161").Value = [double]"-1.305286365105296e-06"
$ws.Range("B162").Value = [double]"0.06162873404420342"
$ws.Range("B163").Value = [double]"0.1919270248274053"
$ws.Range("B164").Value = [double]"-4.184463572873687e-05"
$ws.Range("B165").Value = [double]"-0.02541871493359663"
